$d = $word.ActiveDocument

# Locate the "Marketing Strategy and Data-Driven Insights" paragraph under the
# PARTNER - Siege Analytics entry.
$rng = $d.Content
$found = $rng.Find.Execute("Marketing Strategy and Data-Driven Insights",
                            $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    # New bullet paragraphs to insert immediately after the found paragraph,
    # before the existing "Conducted comprehensive ..." bullet.
    $bullets = @(
        "Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
        "Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
        "Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
    )

    $insertText = ""
    foreach ($bullet in $bullets) {
        $insertText = $insertText + "`r" + [char]0x2022 + " " + $bullet
    }

    # $rng is collapsed to the found text by Find.Execute; inserting text that
    # begins with a paragraph mark after it creates new paragraphs following
    # the "Marketing Strategy and Data-Driven Insights" paragraph.
    $rng.InsertAfter($insertText)
}
